$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The three observation rows (9, 10, 11) had their Id (A), Antal (I),
# Ost (Q) and Nord (R) values cyclically rotated between rows:
#   new row9  = old row11
#   new row10 = old row9
#   new row11 = old row10
# Column "I" (Antal) is stored as text in the source data, so we force
# the cell to text format before writing the value to keep it a string
# instead of Excel auto-converting it to a number.

$ws.Cells.Item(9, 1).Value = 111675585
$ws.Cells.Item(9, 9).NumberFormat = "@"
$ws.Cells.Item(9, 9).Value = "1"
$ws.Cells.Item(9, 17).Value = 690349.9096738817
$ws.Cells.Item(9, 18).Value = 6661440.004307052

$ws.Cells.Item(10, 1).Value = 111675587
$ws.Cells.Item(10, 9).NumberFormat = "@"
$ws.Cells.Item(10, 9).Value = "3"
$ws.Cells.Item(10, 17).Value = 690344.8588249951
$ws.Cells.Item(10, 18).Value = 6661440.743740954

$ws.Cells.Item(11, 1).Value = 111675586
$ws.Cells.Item(11, 9).NumberFormat = "@"
$ws.Cells.Item(11, 9).Value = "2"
$ws.Cells.Item(11, 17).Value = 690348.8581766916
$ws.Cells.Item(11, 18).Value = 6661440.95072202
